# Apply the edits described by the diff to "Todo&Defect.xlsx":
#
#  1) workbook window position: xWindow 3350 -> 4470 (cosmetic; best effort).
#  2) Sheet "Todo ":
#       - row 52: set Scene (col C) to "ui"
#       - new row 53: Issue="loading bar", Description="loading for scene
#         changing", Scene="all", Creator="Fish", PIC="Fish",
#         Create date="02 Apr"
#       - new row 54: Issue="mute at setting and staating page",
#         Scene="all", Creator="Fish", PIC="Fish", Create date="02 Apr"
#       - view: top-left visible cell A28, selection on A14 (best effort).

$wb = $excel.ActiveWorkbook

# --- Workbook window position (best effort; cosmetic window coordinates) ---
$excel.Left   = 4470
$excel.Top    = 460
$excel.Width  = 25600
$excel.Height = 15460

# --- Worksheet edits ---
$ws = $wb.Worksheets.Item("Todo ")
$ws.Activate()

# Existing row 52 gains a Scene value
$ws.Range("C52").Value = "ui"

# New row 53
$ws.Range("A53").Value = "loading bar"
$ws.Range("B53").Value = "loading for scene changing"
$ws.Range("C53").Value = "all"
$ws.Range("D53").Value = "Fish"
$ws.Range("E53").Value = "Fish"
$ws.Range("F53").Value = "02 Apr"

# New row 54 (no Description)
$ws.Range("A54").Value = "mute at setting and staating page"
$ws.Range("C54").Value = "all"
$ws.Range("D54").Value = "Fish"
$ws.Range("E54").Value = "Fish"
$ws.Range("F54").Value = "02 Apr"

# --- View state (best effort) ---
$ws.Range("A28").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A14").Select() | Out-Null
